$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 1. Remove the extra "Sheet1" worksheet (and its Table1 ListObject along
#    with it) - the workbook should end up with only "master-template_type".
# ---------------------------------------------------------------------------
foreach ($s in @($wb.Worksheets)) {
    if ($s.Name -eq "Sheet1") {
        [void]$s.Delete()
    }
}

$ws1 = $wb.Worksheets.Item("master-template_type")

# ---------------------------------------------------------------------------
# 2. Append the new master-data rows (92-121) to the remaining sheet.
# ---------------------------------------------------------------------------
$newRows = @(
    @("RPR_UIN_CARD_TEMPLATE", "UIN card template", "eng"),
    @("RPR_UIN_CARD_TEMPLATE", "قالب بطاقة UIN", "ara"),
    @("RPR_UIN_CARD_TEMPLATE", "Modèle de carte UIN", "fra"),
    @("RPR_UIN_DEAC_SMS", "Template for UIN Deactivation SMS", "eng"),
    @("RPR_UIN_DEAC_SMS", "قالب لتعطيل UIN SMS", "ara"),
    @("RPR_UIN_DEAC_SMS", "Modèle pour SMS de désactivation UIN", "fra"),
    @("RPR_UIN_DEAC_EMAIL", "Template for UIN Deactivation Email", "eng"),
    @("RPR_UIN_DEAC_EMAIL", "قالب لإلغاء تنشيط البريد", "ara"),
    @("RPR_UIN_DEAC_EMAIL", "Modèle pour Email de désactivation UIN", "fra"),
    @("RPR_UIN_REAC_SMS", "Template for UIN Reactivate SMS", "eng"),
    @("RPR_UIN_REAC_SMS", "قالب لـ UIN تنشيط SMS", "ara"),
    @("RPR_UIN_REAC_SMS", "Modèle pour UIN Réactiver SMS", "fra"),
    @("RPR_UIN_REAC_EMAIL", "Template for UIN Reactivate Email", "eng"),
    @("RPR_UIN_REAC_EMAIL", "قالب لـ UIN تنشيط البريد", "ara"),
    @("RPR_UIN_REAC_EMAIL", "Modèle pour UIN Réactiver Email", "fra"),
    @("reg-sms-notification", "Registration Acknowledgement Template", "eng"),
    @("reg-sms-notification", "نموذج شكر التسجيل", "ara"),
    @("reg-sms-notification", "accusé de réception", "fra"),
    @("reg-email-notification", "Registration Acknowledgement Template", "eng"),
    @("reg-email-notification", "نموذج شكر التسجيل", "ara"),
    @("reg-email-notification", "accusé de réception", "fra"),
    @("reg-ack-template-part1", "Registration Acknowledgement Template - Part 1", "eng"),
    @("reg-ack-template-part2", "نموذج شكر التسجيل", "ara"),
    @("reg-ack-template-part3", "accusé de réception", "fra"),
    @("reg-ack-template-part2", "Registration Acknowledgement Template - Part 2", "eng"),
    @("reg-ack-template-part3", "نموذج شكر التسجيل", "ara"),
    @("reg-ack-template-part4", "accusé de réception", "fra"),
    @("reg-ack-template-part3", "Registration Acknowledgement Template - Part 3", "eng"),
    @("reg-ack-template-part4", "نموذج شكر التسجيل", "ara"),
    @("reg-ack-template-part5", "accusé de réception", "fra")
)

$startRow = 92
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $true
    $ws1.Cells.Item($r, 5).Value = "superadmin"
    $ws1.Cells.Item($r, 6).Value = "now()"
}

# ---------------------------------------------------------------------------
# 3. Update the sheet view selection to reflect the new used range, and
#    move the active cell below the last populated row.
# ---------------------------------------------------------------------------
$lastRow = $startRow + $newRows.Length
[void]$ws1.Range("A" + $lastRow + ":XFD1048576").Select()

